# Applies the roster/odds data swap described in the commit:
# "Atualizacao de bases das ligas, do dia: 20-06-2024 as 20:11"
#
# Rows 8-11, 98-99 and 134-135 had their match data (everything except the
# running index in column A, the Div in column C and the Date in column D)
# reshuffled between rows. This script writes the final resulting values for
# every touched cell directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("B8").Value = 5734498
$ws.Range("E8").Value = "Valladolid"
$ws.Range("F8").Value = "Getafe"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = "D"
$ws.Range("L8").Value = 2.1
$ws.Range("M8").Value = 3.4
$ws.Range("N8").Value = 3.4
$ws.Range("O8").Value = 2.05
$ws.Range("P8").Value = 3.3
$ws.Range("Q8").Value = 4
$ws.Range("R8").Value = -0.5
$ws.Range("S8").Value = 2.08
$ws.Range("T8").Value = 1.82
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.82
$ws.Range("W8").Value = 2.08
$ws.Range("X8").Value = -1
$ws.Range("Y8").Value = 2.3
$ws.Range("Z8").Value = -1
$ws.Range("AA8").Value = -1
$ws.Range("AB8").Value = 0.8200000000000001
$ws.Range("AC8").Value = -1
$ws.Range("AD8").Value = 1.08

# Row 9
$ws.Range("B9").Value = 5732768
$ws.Range("E9").Value = "Espanyol"
$ws.Range("F9").Value = "Almeria"
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = "D"
$ws.Range("L9").Value = 3.5
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 1.909
$ws.Range("O9").Value = 2.4
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 2.75
$ws.Range("R9").Value = -0.25
$ws.Range("S9").Value = 2.1
$ws.Range("T9").Value = 1.83
$ws.Range("U9").Value = 3.25
$ws.Range("V9").Value = 2.1
$ws.Range("W9").Value = 1.8
$ws.Range("X9").Value = -1
$ws.Range("Y9").Value = 2.75
$ws.Range("Z9").Value = -1
$ws.Range("AA9").Value = -0.5
$ws.Range("AB9").Value = 0.415
$ws.Range("AC9").Value = 1.1
$ws.Range("AD9").Value = -1

# Row 10
$ws.Range("B10").Value = 5737343
$ws.Range("E10").Value = "Elche"
$ws.Range("F10").Value = "Cadiz"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = "D"
$ws.Range("L10").Value = 2.6
$ws.Range("M10").Value = 3.4
$ws.Range("N10").Value = 2.6
$ws.Range("O10").Value = 2.45
$ws.Range("P10").Value = 3.3
$ws.Range("Q10").Value = 2.875
$ws.Range("R10").Value = -0.25
$ws.Range("S10").Value = 2.05
$ws.Range("T10").Value = 1.75
$ws.Range("U10").Value = 2.5
$ws.Range("V10").Value = 1.91
$ws.Range("W10").Value = 1.99
$ws.Range("X10").Value = -1
$ws.Range("Y10").Value = 2.3
$ws.Range("Z10").Value = -1
$ws.Range("AA10").Value = -0.5
$ws.Range("AB10").Value = 0.375
$ws.Range("AC10").Value = -1
$ws.Range("AD10").Value = 0.99

# Row 11
$ws.Range("B11").Value = 5738382
$ws.Range("E11").Value = "Celta Vigo"
$ws.Range("F11").Value = "Barcelona"
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = "H"
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 3.6
$ws.Range("N11").Value = 2.2
$ws.Range("O11").Value = 2.7
$ws.Range("P11").Value = 3.6
$ws.Range("Q11").Value = 2.45
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 2.06
$ws.Range("T11").Value = 1.84
$ws.Range("U11").Value = 2.5
$ws.Range("V11").Value = 1.95
$ws.Range("W11").Value = 1.95
$ws.Range("X11").Value = 1.7
$ws.Range("Y11").Value = -1
$ws.Range("Z11").Value = -1
$ws.Range("AA11").Value = 1.06
$ws.Range("AB11").Value = -1
$ws.Range("AC11").Value = 0.95
$ws.Range("AD11").Value = -1

# Row 98
$ws.Range("B98").Value = 6809394
$ws.Range("E98").Value = "Celta Vigo"
$ws.Range("F98").Value = "Getafe"
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 2
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = 2
$ws.Range("K98").Value = "D"
$ws.Range("L98").Value = 2.1
$ws.Range("M98").Value = 3.25
$ws.Range("N98").Value = 3.5
$ws.Range("O98").Value = 2
$ws.Range("P98").Value = 3.25
$ws.Range("Q98").Value = 4.2
$ws.Range("R98").Value = -0.5
$ws.Range("S98").Value = 2.01
$ws.Range("T98").Value = 1.89
$ws.Range("U98").Value = 2.25
$ws.Range("V98").Value = 2.07
$ws.Range("W98").Value = 1.83
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = 2.25
$ws.Range("Z98").Value = -1
$ws.Range("AA98").Value = -1
$ws.Range("AB98").Value = 0.8899999999999999
$ws.Range("AC98").Value = 1.07
$ws.Range("AD98").Value = -1

# Row 99
$ws.Range("B99").Value = 6809395
$ws.Range("E99").Value = "CD Alaves"
$ws.Range("F99").Value = "Real Betis"
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 1
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 1
$ws.Range("K99").Value = "D"
$ws.Range("L99").Value = 2.8
$ws.Range("M99").Value = 3.25
$ws.Range("N99").Value = 2.5
$ws.Range("O99").Value = 2.45
$ws.Range("P99").Value = 3.2
$ws.Range("Q99").Value = 3
$ws.Range("R99").Value = -0.25
$ws.Range("S99").Value = 2.11
$ws.Range("T99").Value = 1.79
$ws.Range("U99").Value = 2.25
$ws.Range("V99").Value = 1.99
$ws.Range("W99").Value = 1.91
$ws.Range("X99").Value = -1
$ws.Range("Y99").Value = 2.2
$ws.Range("Z99").Value = -1
$ws.Range("AA99").Value = -0.5
$ws.Range("AB99").Value = 0.395
$ws.Range("AC99").Value = -0.5
$ws.Range("AD99").Value = 0.455

# Row 134
$ws.Range("B134").Value = 6809255
$ws.Range("E134").Value = "Osasuna"
$ws.Range("F134").Value = "Las Palmas"
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 1
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = "D"
$ws.Range("L134").Value = 1.75
$ws.Range("M134").Value = 3.6
$ws.Range("N134").Value = 4.75
$ws.Range("O134").Value = 1.909
$ws.Range("P134").Value = 3.2
$ws.Range("Q134").Value = 4.5
$ws.Range("R134").Value = -0.5
$ws.Range("S134").Value = 1.97
$ws.Range("T134").Value = 1.93
$ws.Range("U134").Value = 2
$ws.Range("V134").Value = 1.75
$ws.Range("W134").Value = 2.05
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = 2.2
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = -1
$ws.Range("AB134").Value = 0.9299999999999999
$ws.Range("AC134").Value = 0
$ws.Range("AD134").Value = 0

# Row 135
$ws.Range("B135").Value = 6809414
$ws.Range("E135").Value = "Granada"
$ws.Range("F135").Value = "Getafe"
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 1
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = 1
$ws.Range("K135").Value = "D"
$ws.Range("L135").Value = 2.5
$ws.Range("M135").Value = 3.1
$ws.Range("N135").Value = 3
$ws.Range("O135").Value = 2.45
$ws.Range("P135").Value = 3.1
$ws.Range("Q135").Value = 3.1
$ws.Range("R135").Value = -0.25
$ws.Range("S135").Value = 2.11
$ws.Range("T135").Value = 1.79
$ws.Range("U135").Value = 2.25
$ws.Range("V135").Value = 2.02
$ws.Range("W135").Value = 1.88
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 2.1
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = -0.5
$ws.Range("AB135").Value = 0.395
$ws.Range("AC135").Value = -0.5
$ws.Range("AD135").Value = 0.4399999999999999
